$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# 1) "Il s'adapte a la vitesse du processeur" paragraph (section 5 -
#    FONCTIONNEMENT): the sentence about the game loop timing gets two new
#    clauses appended - "en conséquence" tacked onto the existing sentence,
#    and a brand-new sentence explaining the consequence.
# ---------------------------------------------------------------------------
$d.Content.Find.Execute(
    "les déplacements et rotations des personnages du jeu.",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "les déplacements et rotations des personnages du jeu en conséquence. Ainsi, les mouvements des personnages restent uniformes même en cas de ralentissement du processeur.",
    2)

# ---------------------------------------------------------------------------
# 2) Final "Bilan personnel" paragraph: "je travaillais sur un projet" ->
#    "je travaillais en groupe sur un projet".
# ---------------------------------------------------------------------------
$d.Content.Find.Execute(
    "je travaillais sur un projet en langage Python.",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "je travaillais en groupe sur un projet en langage Python.",
    2)

# ---------------------------------------------------------------------------
# 3) The hidden "_GoBack" bookmark (marking Word's last-edit position) moves
#    from the empty paragraph just before "7 - PROLONGEMENT POSSIBLE" to the
#    very end of the document (end of the last "Bilan personnel" paragraph),
#    matching where the newest edit above took place.
#
#    Adding a bookmark named "_GoBack" automatically removes the previous
#    one elsewhere in the document (Word only ever keeps a single "_GoBack"
#    bookmark). A zero-length bookmark placed exactly at the last character
#    boundary of a paragraph can get mis-anchored, so we briefly insert a
#    sentinel character, bookmark just before it, then delete the sentinel -
#    leaving the bookmark correctly collapsed at the true end of the text.
# ---------------------------------------------------------------------------
$lastPara = $d.Paragraphs.Item($d.Paragraphs.Count)
$endRng = $lastPara.Range.Duplicate
$endRng.MoveEnd(1, -1)
$endRng.Collapse(0)
$endRng.InsertAfter("X")

$bmRng = $d.Range($endRng.Start, $endRng.Start)
$d.Bookmarks.Add("_GoBack", $bmRng)

$sentinel = $d.Range($endRng.Start, $endRng.Start + 1)
$sentinel.Text = ""
